$d = $word.ActiveDocument
$newText = "Fechas de la campaña para Perseo: Perseus: 16-25 de enero, del 7 al 16 de noviembre, del 6 al 15 de diciembre"

$targets = @(3, 53, 86, 121)
foreach ($idx in $targets) {
    $p = $d.Paragraphs($idx)
    $rng = $p.Range
    $rng.MoveEnd(1, -1)
    $rng.Delete()
    $rng2 = $p.Range
    $rng2.MoveEnd(1, -1)
    $rng2.InsertAfter($newText)
}
